# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets
# to match refreshed data output (regenerated gh-pages output).

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1158
$ws1.Range("F7").Value  = 12219
$ws1.Range("F8").Value  = 55
$ws1.Range("F11").Value = 12018
$ws1.Range("F12").Value = 4799
$ws1.Range("F13").Value = 2630
$ws1.Range("F15").Value = 48

# Sheet "全部类型": row -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1158
$ws4.Range("F9").Value  = 12219
$ws4.Range("F10").Value = 55
$ws4.Range("F13").Value = 12018
$ws4.Range("F14").Value = 4799
$ws4.Range("F15").Value = 2630
$ws4.Range("F17").Value = 48
